$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.753.48'
$ws.Range('E2').Value = '  -7.25%  '
$ws.Range('D3').Value = '2.543.99'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '296.57'
$ws.Range('E5').Value = '  -4.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.39'
$ws.Range('E6').Value = '  -5.63%  '
$ws.Range('E7').Value = '  -4.02%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').Value = '  -5.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.80'
$ws.Range('E10').Value = '  -8.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  -4.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.67'
$ws.Range('E12').Value = '  -5.42%  '
$ws.Range('D13').Value = '2.925.85'
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.106'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = '2.537.42'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.865'
$ws.Range('E16').Value = '  -5.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.09'
$ws.Range('E17').Value = '  -4.83%  '
$ws.Range('D18').Value = '42.751.82'
$ws.Range('E18').Value = '  -7.49%  '
$ws.Range('D19').Value = '0.0₃0966'
$ws.Range('E19').Value = '  -4.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.55'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.39'
$ws.Range('E21').Value = '  -3.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.98'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '258.23'
$ws.Range('E23').Value = '  -5.45%  '
$ws.Range('E24').Value = '  -5.25%  '
$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.59'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.12'
$ws.Range('E26').Value = '  -3.21%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('E29').Value = '  -7.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.05'
$ws.Range('E30').Value = '  -5.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.88'
$ws.Range('E31').Value = '  -5.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.97'
$ws.Range('E32').Value = '  -3.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.39'
$ws.Range('E33').Value = '  -6.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.15'
$ws.Range('E34').Value = '  -3.44%  '
$ws.Range('E35').Value = '  -3.46%  '
$ws.Range('E36').Value = '  -5.26%  '
$ws.Range('E37').Value = '  -8.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.16'
$ws.Range('E38').Value = '  +4.97%  '
$ws.Range('E39').Value = '  -3.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '15.86'
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.42'
$ws.Range('E41').Value = '  -4.87%  '
$ws.Range('E42').Value = '  -7.10%  '
$ws.Range('D43').Value = '2.076.02'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('E44').Value = '  -4.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '84.69'
$ws.Range('E46').Value = '  -11.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.58'
$ws.Range('E47').Value = '  +3.13%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.70'
$ws.Range('E48').Value = '  -3.12%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.75'
$ws.Range('E49').Value = '  -9.39%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.782.98'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.80'
$ws.Range('E51').Value = '  -5.47%  '
